$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 39 (the "Media" row), shifting everything below down by one.
$ws.Rows.Item(39).Insert()

# Populate the new row 39 with the "Synced Lyrics" tag mapping.
$ws.Cells.Item(39, 1).Value = "Synced Lyrics"
$ws.Cells.Item(39, 2).Value = "syncedlyrics:language:description"
$ws.Cells.Item(39, 3).Value = "SYLT:description"
$ws.Cells.Item(39, 4).Value = "n/a"
$ws.Cells.Item(39, 5).Value = "n/a"
$ws.Cells.Item(39, 6).Value = "n/a"
$ws.Cells.Item(39, 7).Value = "n/a"
$ws.Cells.Item(39, 8).Value = "n/a"

# Match the formatting of the surrounding data rows (column A uses a bold/bordered
# style, columns B-H use a wrapped/top-aligned style) by copying formats down from
# the row directly above (row 38, "Lyrics [4]").
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)

$ws.Range("B38:H38").Copy()
$ws.Range("B39:H39").PasteSpecial(-4122)

$excel.CutCopyMode = 0
